# Update mods data [2026-02-05 15:30:38]
# Appends a new row (row 87) to the ModCounts sheet with the latest mod count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 86
$newRow = $lastRow + 1

# Force column A to be treated as literal text (not auto-parsed as a date)
# while we set the new cell values.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2026/02/05"
$ws.Range("B$newRow").Value = "逃离鸭科夫"
$ws.Range("C$newRow").Value = 1174

# Re-apply the same look/format used by the rest of the data rows
# (copy only formatting from the previous row, leaving the values untouched).
$ws.Range("A$lastRow`:C$lastRow").Copy()
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122)
